$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "1.03" label in A9 -> "01.03" ---
# A plain Value assignment of a numeric-looking string like "01.03" gets
# auto-coerced into the number 1.03 (losing the leading zero / text type),
# so we stage the literal text in a scratch cell (quote-prefixed so Excel
# keeps it as text), copy only its value into the target cell, then clean
# the scratch cell back up.
$ws.Range("K1").Value = "'01.03"
$ws.Range("K1").Copy()
$ws.Range("A9").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("K1").Clear()

# --- Add the new row 10 for the 02.03 data ---
# First clone A9's formatting (border/font/alignment) down into A10 so the
# new label cell matches the other date-label cells...
$ws.Range("A9").Copy($ws.Range("A10"))
# ...then overwrite just the value with "02.03" using the same
# text-safe staging trick as above.
$ws.Range("K1").Value = "'02.03"
$ws.Range("K1").Copy()
$ws.Range("A10").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("K1").Clear()

$ws.Range("B10").Value = 742
$ws.Range("C10").Value = 166
$ws.Range("D10").Value = 927
$ws.Range("E10").Value = 1835
$ws.Range("F10").Value = 149
$ws.Range("G10").Value = 52
$ws.Range("H10").Value = 2036
$ws.Range("I10").Value = 23345
